$d = $word.ActiveDocument

$bodyXml = @'
<w:p w:rsidR="00652232" w:rsidRDefault="00652232" w:rsidP="00652232"><w:r><w:t>Week 3:</w:t></w:r></w:p><w:p w:rsidR="00652232" w:rsidRDefault="00652232" w:rsidP="00652232"><w:r><w:t xml:space="preserve">1a: Dit komt voor bij de methodes </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>aantalArtikelen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">() en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hoeveelheidGeldInKassa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>().</w:t></w:r><w:r><w:tab/></w:r></w:p><w:p w:rsidR="00AB5F13" w:rsidRDefault="00D55D22" w:rsidP="00652232"><w:r><w:t>2a</w:t></w:r><w:r w:rsidR="00AB5F13"><w:t xml:space="preserve">: </w:t></w:r><w:r w:rsidR="00652232"><w:t xml:space="preserve">Het is handig om de methodes private te maken omdat dan alleen de klasse waarin die methode gemaakt is de gegevens zal gebruiken. Wanneer je dus iets verandert aan de methodes, zal dit geen invloed hebben op andere </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00652232"><w:t>klasses</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00652232"><w:t xml:space="preserve">. Dit maakt de code een stuk overzichtelijker. </w:t></w:r></w:p><w:p w:rsidR="00AB5F13" w:rsidRDefault="00652232" w:rsidP="00652232"><w:r><w:t xml:space="preserve">2b: In een </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>HashMap</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> heb je </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>key's</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>value's</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. In </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>HashSet</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> heb je alleen </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>value's</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p><w:p w:rsidR="00AB5F13" w:rsidRDefault="00652232" w:rsidP="00652232"><w:r><w:t xml:space="preserve">3a: </w:t></w:r><w:r><w:t xml:space="preserve">De </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>constructor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> maakt als eerst een nieuwe Kantine en een nieuwe Random aan. Hierna wordt er een array gemaakt met AANTAL_ARTIKELEN waardes die allemaal tussen MIN_ARTIKELEN_PER_SOORT en MAX_ARTIKELEN_PER_SOORT liggen.</w:t></w:r><w:r><w:t xml:space="preserve"> Hierna wordt er een </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>KantineAanbod</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> gemaakt waarbij met de artikelnamen, artikelprijzen en hoeveelheden meegeeft. Tot slot wordt het kantineaanbod </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>geset</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> voor een Kantine</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p w:rsidR="00652232" w:rsidRDefault="00652232" w:rsidP="00652232"><w:r><w:t>3b: Aangezien er +min achter staat, wordt de waarde opgehoogd met het minimum. Normaal is de minimu</w:t></w:r><w:r w:rsidR="00D55D22"><w:t xml:space="preserve">m waarde 0. Dus nu is de minimum </w:t></w:r><w:r><w:t>waarde 0 + min = min.</w:t></w:r></w:p><w:p w:rsidR="00652232" w:rsidRDefault="00652232" w:rsidP="00D55D22"><w:r><w:t>Het maximum komt tussen haakjes te staan. Maar omdat er later nog weer de minimum waarde bij opgeteld wordt, moet de minimum</w:t></w:r><w:r w:rsidR="00D55D22"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>waarde hier weer van het maximum afgehaald worden.</w:t></w:r></w:p><w:p w:rsidR="002200E5" w:rsidRDefault="00652232" w:rsidP="00652232"><w:r><w:t xml:space="preserve">Normaal zou </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nextInt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> TOT de ingevoerde </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>value</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> gaan. Maar in dit geval gaat hij tot en met de ingevoerde </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>value</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> gaan. Door de +1.</w:t></w:r></w:p><w:p w:rsidR="00D55D22" w:rsidRDefault="00D55D22" w:rsidP="00652232"/><w:p w:rsidR="00D55D22" w:rsidRDefault="00D55D22" w:rsidP="00652232"><w:r><w:t>Week 4:</w:t></w:r></w:p><w:p w:rsidR="00D55D22" w:rsidRDefault="00D55D22" w:rsidP="00652232"><w:r><w:t xml:space="preserve">1c: Er hoeven geen waarden mee te worden gegeven aan het object. Een </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>constructor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is in dit geval dus niet nodig. Het enige wat we nodig hebben zijn de methodes.</w:t></w:r></w:p><w:p w:rsidR="00871EF1" w:rsidRDefault="00871EF1" w:rsidP="00652232"><w:r><w:t xml:space="preserve">1d: Deze 2 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>methode’s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> voeren niks uit op een instantie van de klasse Administratie.</w:t></w:r></w:p><w:p w:rsidR="00871EF1" w:rsidRDefault="00871EF1" w:rsidP="00652232"><w:r><w:t xml:space="preserve">1e: Nu de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>constructor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> private is kan men niet meer van ‘buiten af’ een nieuwe instantie van de klasse Administratie aanmaken.</w:t></w:r></w:p><w:p w:rsidR="00E319C6" w:rsidRDefault="00E319C6" w:rsidP="00652232"><w:r><w:t>1</w:t></w:r><w:r w:rsidR="00804B60"><w:t xml:space="preserve">i: Als het niet </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00804B60"><w:t>static</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00804B60"><w:t xml:space="preserve"> is, hoeft het ook niet overal hetzelfde te zijn. Dus kan hij niet </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00804B60"><w:t>final</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00804B60"><w:t xml:space="preserve"> worden omdat het eventueel nog gewijzigd zou kunnen worden.</w:t></w:r></w:p><w:p w:rsidR="00804B60" w:rsidRDefault="00804B60" w:rsidP="00652232"><w:r><w:t>1j: Nu zou het nog weer aangepast kun</w:t></w:r><w:r w:rsidR="008837A4"><w:t>nen worden, wat w</w:t></w:r><w:r><w:t>e</w:t></w:r><w:r w:rsidR="008837A4"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>niet willen.</w:t></w:r></w:p>
'@

$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$d.Content.InsertXML($xmlFrag)

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
